$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A424").NumberFormat = "@"
$ws.Range("A424").Value = '14506209'
$ws.Range("A424").ClearFormats()
$ws.Range("B424").NumberFormat = "@"
$ws.Range("B424").Value = '2025-08-25'
$ws.Range("B424").ClearFormats()
$ws.Range("C424").Value = 'Cameron Norrie'
$ws.Range("D424").Value = 'Sebastian Korda'
$ws.Range("E424").Value = 'Gana Cameron Norrie'
$ws.Range("F424").Value = 2.63

$ws.Range("A425").NumberFormat = "@"
$ws.Range("A425").Value = '14520001'
$ws.Range("A425").ClearFormats()
$ws.Range("B425").NumberFormat = "@"
$ws.Range("B425").Value = '2025-08-25'
$ws.Range("B425").ClearFormats()
$ws.Range("C425").Value = 'Sebastián Báez'
$ws.Range("D425").Value = 'Lloyd Harris'
$ws.Range("E425").Value = 'Gana Sebastián Báez'
$ws.Range("F425").Value = 2.5

$ws.Range("A426").NumberFormat = "@"
$ws.Range("A426").Value = '14506218'
$ws.Range("A426").ClearFormats()
$ws.Range("B426").NumberFormat = "@"
$ws.Range("B426").Value = '2025-08-25'
$ws.Range("B426").ClearFormats()
$ws.Range("C426").Value = 'Ugo Humbert'
$ws.Range("D426").Value = 'Adam Walton'
$ws.Range("E426").Value = 'Gana Ugo Humbert'
$ws.Range("F426").Value = 1.4

$ws.Range("A427").NumberFormat = "@"
$ws.Range("A427").Value = '14506188'
$ws.Range("A427").ClearFormats()
$ws.Range("B427").NumberFormat = "@"
$ws.Range("B427").Value = '2025-08-25'
$ws.Range("B427").ClearFormats()
$ws.Range("C427").Value = 'Mattia Bellucci'
$ws.Range("D427").Value = 'Juncheng Shang'
$ws.Range("E427").Value = 'Gana Juncheng Shang'
$ws.Range("F427").Value = 2.2

$ws.Range("A428").NumberFormat = "@"
$ws.Range("A428").Value = '14520002'
$ws.Range("A428").ClearFormats()
$ws.Range("B428").NumberFormat = "@"
$ws.Range("B428").Value = '2025-08-25'
$ws.Range("B428").ClearFormats()
$ws.Range("C428").Value = 'Jan-Lennard Struff'
$ws.Range("D428").Value = 'Mackenzie McDonald'
$ws.Range("E428").Value = 'Gana Jan-Lennard Struff'
$ws.Range("F428").Value = 1.91

$ws.Range("A429").NumberFormat = "@"
$ws.Range("A429").Value = '14520015'
$ws.Range("A429").ClearFormats()
$ws.Range("B429").NumberFormat = "@"
$ws.Range("B429").Value = '2025-08-25'
$ws.Range("B429").ClearFormats()
$ws.Range("C429").Value = 'Jaume Munar'
$ws.Range("D429").Value = 'Jaime Faria'
$ws.Range("E429").Value = 'Gana Jaime Faria'
$ws.Range("F429").Value = 3.5

$ws.Range("A430").NumberFormat = "@"
$ws.Range("A430").Value = '14506224'
$ws.Range("A430").ClearFormats()
$ws.Range("B430").NumberFormat = "@"
$ws.Range("B430").Value = '2025-08-25'
$ws.Range("B430").ClearFormats()
$ws.Range("C430").Value = 'Jenson Brooksby'
$ws.Range("D430").Value = 'Aleksandar Vukic'
$ws.Range("E430").Value = 'Gana Aleksandar Vukic'
$ws.Range("F430").Value = 3.2

$ws.Range("A431").NumberFormat = "@"
$ws.Range("A431").Value = '14506223'
$ws.Range("A431").ClearFormats()
$ws.Range("B431").NumberFormat = "@"
$ws.Range("B431").Value = '2025-08-25'
$ws.Range("B431").ClearFormats()
$ws.Range("C431").Value = 'Gabriel Diallo'
$ws.Range("D431").Value = 'Damir Dzumhur'
$ws.Range("E431").Value = 'Gana Damir Dzumhur'
$ws.Range("F431").Value = 4

$ws.Range("A432").NumberFormat = "@"
$ws.Range("A432").Value = '14506195'
$ws.Range("A432").ClearFormats()
$ws.Range("B432").NumberFormat = "@"
$ws.Range("B432").Value = '2025-08-25'
$ws.Range("B432").ClearFormats()
$ws.Range("C432").Value = 'Sebastian Ofner'
$ws.Range("D432").Value = 'Casper Ruud'
$ws.Range("E432").Value = 'Gana Sebastian Ofner'
$ws.Range("F432").Value = 4.33

$ws.Range("A433").NumberFormat = "@"
$ws.Range("A433").Value = '14519801'
$ws.Range("A433").ClearFormats()
$ws.Range("B433").NumberFormat = "@"
$ws.Range("B433").Value = '2025-08-25'
$ws.Range("B433").ClearFormats()
$ws.Range("C433").Value = 'Leolia Jeanjean'
$ws.Range("D433").Value = 'Priscilla Hon'
$ws.Range("E433").Value = 'Gana Priscilla Hon'
$ws.Range("F433").Value = 1.91

$ws.Range("A434").NumberFormat = "@"
$ws.Range("A434").Value = '14506264'
$ws.Range("A434").ClearFormats()
$ws.Range("B434").NumberFormat = "@"
$ws.Range("B434").Value = '2025-08-25'
$ws.Range("B434").ClearFormats()
$ws.Range("C434").Value = 'Madison Keys'
$ws.Range("D434").Value = 'Renata Zarazua'
$ws.Range("E434").Value = 'Gana Renata Zarazua'
$ws.Range("F434").Value = 9

$ws.Range("A435").NumberFormat = "@"
$ws.Range("A435").Value = '14506239'
$ws.Range("A435").ClearFormats()
$ws.Range("B435").NumberFormat = "@"
$ws.Range("B435").Value = '2025-08-25'
$ws.Range("B435").ClearFormats()
$ws.Range("C435").Value = 'Anna Kalinskaya'
$ws.Range("D435").Value = 'Clervie Ngounoue'
$ws.Range("E435").Value = 'Gana Clervie Ngounoue'
$ws.Range("F435").Value = 3.75

$ws.Range("A436").NumberFormat = "@"
$ws.Range("A436").Value = '14506267'
$ws.Range("A436").ClearFormats()
$ws.Range("B436").NumberFormat = "@"
$ws.Range("B436").Value = '2025-08-25'
$ws.Range("B436").ClearFormats()
$ws.Range("C436").Value = 'Taylor Townsend'
$ws.Range("D436").Value = 'Antonia Ružić'
$ws.Range("E436").Value = 'Gana Antonia Ružić'
$ws.Range("F436").Value = 2.3

$ws.Range("A437").NumberFormat = "@"
$ws.Range("A437").Value = '14519997'
$ws.Range("A437").ClearFormats()
$ws.Range("B437").NumberFormat = "@"
$ws.Range("B437").Value = '2025-08-25'
$ws.Range("B437").ClearFormats()
$ws.Range("C437").Value = 'Darja Semenistaja'
$ws.Range("D437").Value = 'Peyton Stearns'
$ws.Range("E437").Value = 'Gana Darja Semenistaja'
$ws.Range("F437").Value = 3.5

$ws.Range("A438").NumberFormat = "@"
$ws.Range("A438").Value = '14528381'
$ws.Range("A438").ClearFormats()
$ws.Range("B438").NumberFormat = "@"
$ws.Range("B438").Value = '2025-08-26'
$ws.Range("B438").ClearFormats()
$ws.Range("C438").Value = 'Nicolai Budkov Kjaer'
$ws.Range("D438").Value = 'Cedrik-Marcel Stebe'
$ws.Range("E438").Value = 'Gana Cedrik-Marcel Stebe'
$ws.Range("F438").Value = 3.4

$ws.Range("A439").NumberFormat = "@"
$ws.Range("A439").Value = '14534329'
$ws.Range("A439").ClearFormats()
$ws.Range("B439").NumberFormat = "@"
$ws.Range("B439").Value = '2025-08-26'
$ws.Range("B439").ClearFormats()
$ws.Range("C439").Value = 'Valentin Vacherot'
$ws.Range("D439").Value = 'Martin Landaluce'
$ws.Range("E439").Value = 'Gana Valentin Vacherot'
$ws.Range("F439").Value = 2.2

$ws.Range("A440").NumberFormat = "@"
$ws.Range("A440").Value = '14527268'
$ws.Range("A440").ClearFormats()
$ws.Range("B440").NumberFormat = "@"
$ws.Range("B440").Value = '2025-08-26'
$ws.Range("B440").ClearFormats()
$ws.Range("C440").Value = 'Masamichi Imamura'
$ws.Range("D440").Value = 'Oliver Crawford'
$ws.Range("E440").Value = 'Gana Masamichi Imamura'
$ws.Range("F440").Value = 3.75

$ws.Range("A441").NumberFormat = "@"
$ws.Range("A441").Value = '14528388'
$ws.Range("A441").ClearFormats()
$ws.Range("B441").NumberFormat = "@"
$ws.Range("B441").Value = '2025-08-25'
$ws.Range("B441").ClearFormats()
$ws.Range("C441").Value = 'Gilles Arnaud Bailly'
$ws.Range("D441").Value = 'Maxim Mrva'
$ws.Range("E441").Value = 'Gana Gilles Arnaud Bailly'
$ws.Range("F441").Value = 1.83
